$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.674.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.597.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("E6").Value = "  +1.34%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("E10").Value = "  +0.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0843"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.822.62"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.570.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("E14").Value = "  -0.80%  "
$ws.Range("E15").Value = "  -1.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.663.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "208.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  -3.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  -2.76%  "
$ws.Range("E28").Value = "  +1.81%  "
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("E32").Value = "  -0.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.657"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.295.70"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.79%  "
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.842"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.34%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  +1.67%  "
$ws.Range("E42").Value = "  +1.01%  "
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.734.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.898"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.18%  "
$ws.Range("E48").Value = "  -0.50%  "
$ws.Range("E49").Value = "  +1.66%  "
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.09%  "
